$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 490, shifting existing rows 490-559 down to 491-560
$ws.Rows.Item(490).Insert()

# Populate the new row 490 with fresh data
$ws.Cells.Item(490, 1).Value = 7
$ws.Cells.Item(490, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(490, 3).Value = "Ñuble"
$ws.Cells.Item(490, 4).Value = 45131
$ws.Cells.Item(490, 5).Value = 16
$ws.Cells.Item(490, 6).Value = 100112023
$ws.Cells.Item(490, 7).Value = "Brócoli"
$ws.Cells.Item(490, 8).Value = "Sin especificar"
$ws.Cells.Item(490, 9).Value = "Primera"
$ws.Cells.Item(490, 10).Value = 250
$ws.Cells.Item(490, 11).Value = 1000
$ws.Cells.Item(490, 12).Value = 1000
$ws.Cells.Item(490, 13).Value = 1000
$ws.Cells.Item(490, 14).Value = "$/unidad"
$ws.Cells.Item(490, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(490, 16).Value = 1000
$ws.Cells.Item(490, 17).Value = 1
$ws.Cells.Item(490, 18).Value = "Hortaliza"
